# sheet_with_same_numeric_value_date_formatted_differently.xlsx
# Commit: "Add DateInterval support in XLSX Reader. Support negative durations in XLSX writer."
#
# The meaningful edit is: cell I3 (row 3) used the built-in elapsed-time
# format "[h]:mm:ss" (numFmtId 46). It is changed to the plain (non-elapsed)
# time format "h:mm:ss" so the reader/writer test fixture exercises the
# difference between elapsed-duration formats and ordinary time formats
# used for DateInterval support.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change I3's number format from the elapsed-time "[h]:mm:ss" to "h:mm:ss".
$ws.Range("I3").NumberFormat = "h:mm:ss"

# The author's selection was on I3 when the file was last saved.
$ws.Range("I3").Select()
